$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.599.86"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "3.399.37"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.14"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.68"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.400.89"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").Value = "  +8.03%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +4.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "3.990.73"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000190"
$ws.Range("E15").Value = "  +5.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.45"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "63.617.17"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "3.387.23"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.11"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.72"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.07"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.97"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.535"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  +22.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  +9.67%  "
$ws.Range("E31").Value = "  +6.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.32"
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.86"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.31"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0759"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.903.80"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.70"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0316"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.18"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.37"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.754"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.23"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").Value = "  +21.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.45"
$ws.Range("E51").Value = "  +2.57%  "
